$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.285500000000004
$ws.Range("D3").Value = -7.250199999999996
$ws.Range("A12").Value = -21.60250000000001
$ws.Range("B14").Value = 7.176599999999994
$ws.Range("D20").Value = -7.505400000000004
$ws.Range("D25").Value = -8.0365
$ws.Range("B26").Value = 4.261300000000007
$ws.Range("A27").Value = -21.77439999999999
$ws.Range("D30").Value = -7.5851
$ws.Range("B31").Value = 5.162600000000005
$ws.Range("A32").Value = -21.40580000000001
$ws.Range("B35").Value = 9.182399999999996
$ws.Range("A36").Value = -19.9242
$ws.Range("B37").Value = 8.881600000000002
$ws.Range("A38").Value = -19.4021
$ws.Range("D44").Value = -7.221800000000004
$ws.Range("B45").Value = 6.467699999999999
$ws.Range("A46").Value = -21.48089999999999
$ws.Range("D47").Value = -7.5076
$ws.Range("B52").Value = 5.107600000000002
$ws.Range("A54").Value = -21.7016
$ws.Range("A55").Value = -22.37660000000001
$ws.Range("A56").Value = -22.28700000000002
$ws.Range("B57").Value = 4.880499999999994
$ws.Range("D58").Value = -8.171799999999996
$ws.Range("A67").Value = -21.43499999999998
$ws.Range("A69").Value = -21.54589999999997
$ws.Range("A72").Value = -21.95650000000001
$ws.Range("D78").Value = -7.587700000000004
$ws.Range("B81").Value = 5.9973
$ws.Range("A83").Value = -21.5605
$ws.Range("B83").Value = 5.691800000000005
$ws.Range("D84").Value = -8.568000000000005
$ws.Range("A86").Value = -22.02830000000001
$ws.Range("D89").Value = -6.210299999999997
$ws.Range("A91").Value = -21.35450000000002
$ws.Range("D91").Value = -6.162199999999998
$ws.Range("D92").Value = -6.103199999999998
$ws.Range("A93").Value = -21.2635
$ws.Range("D96").Value = -7.4321
$ws.Range("A99").Value = -20.23209999999999
$ws.Range("B100").Value = 5.332899999999999
$ws.Range("B102").Value = 8.922300000000002
$ws.Range("D102").Value = -8.046900000000001
